$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Heading "#6 [star] Last Year's Goals" -> append
#    " Progress (After Prep Tasks)" as two new runs: a plain space run,
#    then a styled run using a new "normaltextrun" character style.
# ---------------------------------------------------------------------

# Register the new character style used by the appended run.
$normalTextRunStyle = $d.Styles.Add("normaltextrun", 2)
$normalTextRunStyle.BaseStyle = "DefaultParagraphFont"

$goalsRange = $d.Content
$goalsRange.Find.Execute("Last Year" + [char]8217 + "s Goals") | Out-Null
$goalsRange.Expand(4) | Out-Null  ## wdParagraph - grab the whole heading paragraph

# Insert a plain " " run right before the paragraph mark.
$insertPos = $goalsRange.End - 1
$spaceRange = $d.Range($insertPos, $insertPos)
$spaceRange.InsertAfter(" ")

# Insert the styled "Progress (After Prep Tasks)" run right after the space.
$insertPos2 = $insertPos + 1
$textRange = $d.Range($insertPos2, $insertPos2)
$textRange.InsertAfter("Progress (After Prep Tasks)")

# Apply the character style + font to just the text we inserted.
$progressRange = $d.Range($insertPos2, $insertPos2 + "Progress (After Prep Tasks)".Length)
$progressRange.Style = "normaltextrun"
$progressRange.Font.Name = "Aptos Display"

# ---------------------------------------------------------------------
# 2. Heading "#1" + "8" + " " + "Client Love" -> collapse to one run
#    "#18 Client Love".
# ---------------------------------------------------------------------

$clientLoveFind = $d.Content
$clientLoveFind.Find.Execute("Client Love") | Out-Null
$clientLoveFind.Expand(4) | Out-Null  ## wdParagraph

$headingRange = $d.Range($clientLoveFind.Start, $clientLoveFind.End - 1)
$headingText = $headingRange.Text
$headingRange.Text = ""
$headingRange2 = $d.Range($clientLoveFind.Start, $clientLoveFind.Start)
$headingRange2.Text = $headingText

# ---------------------------------------------------------------------
# 3. Footer "Last saved" date/time text update.
# ---------------------------------------------------------------------

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $footers = $d.Sections.Item($s).Footers
    for ($i = 1; $i -le $footers.Count; $i++) {
        $footer = $footers.Item($i)
        if ($footer.Exists) {
            $footer.Range.Find.Execute("22/04/2024 11:41 AM", $true, $false, $false, $false, $false, `
                $true, 1, $false, "24/04/2024 12:43 PM", 2) | Out-Null
        }
    }
}
